$d = $word.ActiveDocument

# --- Paragraph 1: date line + paper title (joined by a line break) ---
$d.Paragraphs.Item(1).Range.Text = 'המאמר היומי של מייק: 14.06.25' + [char]11 + 'Is Stochastic Gradient Descent Effective? A PDE Perspective on Machine Learning Processes'

# --- Paragraph 2 ---
$d.Paragraphs.Item(2).Range.Text = 'המאמר הוא די כבד אבל ניסיתי להנגיש את הסקירה כך שתהיה מובנת (גם אני לא צללתי עמוק מדי שם - המאמר באמת מורכב).'

# --- Paragraph 3 ---
$d.Paragraphs.Item(3).Range.Text = 'יש משהו מתעתע בפשטות של Stochastic Gradient Descent או SGD בקצרה. כבר שנים שהוא הליבה של למידת מכונה(ML), ובמיוחד של למידה עמוקה, אבל התשובות לשאלה למה הוא בעצם עובד נותרו בגדר אינטואיציה לא מספקת. נדמה שכל ניסיון להסביר את ההצלחה של SGD חוזר בסופו של דבר לאמירות מעורפלות כמו "הוא מוצא מינימות שטוחות" או "הרעש עוזר לצאת ממינימום מקומי". המאמר שאני סוקר היום, מנסה לעשות סדר ובשונה מרוב העבודות בתחום, הוא מציע זווית חדשה לגמרי: הוא מתאר את SGD כתהליך דיפוזיוני שמתפתח בזמן, דרך עדשה של משוואות דיפרנציאליות חלקיות (PDEs).'
$d.Content.Find.Execute('(PDEs).', $true, $false, $false, $false, $false, $true, 1, $false, '(PDEs).', 2) | Out-Null

# --- Paragraph 4 ---
$d.Paragraphs.Item(4).Range.Text = 'המחברים מבקשים לשנות את הדרך שבה אנחנו מבינים את הדינמיקה של למידה. לא עוד מעקב אחרי נקודה במרחב המשקולות שמתגלגלת בתוך משטח לוס (loss landscape), אלא תיאור מלא של ההתפלגות ההסתברותית של כל האפשרויות כלומר צפיפות (במהלך תהליך הלמידה) על פני המרחב, שמתפתחת בזמן. אם אתם מגיעים מתחום הפיזיקה המתמטית, זה יזכיר לכם מיד את משוואת פוקר־פלאנק, שמתארת איך חלקיקים נעים ונפזרים במערכת נתונה. הרעיון כאן הוא דומה: המשקולות הם כמו חלקיקים, והם נעות על פי הגרדיאנט של פונקציית לוס, עם קצת רעש שנובע מהאופיין הסטוכסטי בו(בחירת מיני-באצ''ים) של SGD.'
$d.Content.Find.Execute(') של SGD.', $true, $false, $false, $false, $false, $true, 1, $false, ') של SGD.', 2) | Out-Null

# --- Paragraph 5 ---
$d.Paragraphs.Item(5).Range.Text = 'מה שמעניין הוא שהמודל הפיזיקלי הזה הוא לא רק שהוא מחקה את מה ש-SGD עושה, אלא מראה מדוע הוא מצליח. למשל, כאשר מסתכלים על האנרגיה הקינטית של המערכת, רואים שהרעש האקראי שנובע מהסטוכסטיות של הבחירה במיני־באטצ''ים לא סתם "מוסיף רעש" אלא משחק תפקיד קריטי ביציבות: הוא מאזן את ההתקדמות כך שלא נגלוש מהר מדי או ניתקע במקומות לא יציבים. המחברים ממש מראים כיצד יש מגבלות אנרגטיות שמכתיבות את הקצב שבו אפשר ללמוד, וקושרות בין כמות הרעש לבין עומק הירידה באובדן.'

# --- Paragraph 6 ---
$d.Paragraphs.Item(6).Range.Text = 'יש כאן גם הבחנה מושגית חדה בין שתי גישות להבנת תהליכי למידה: הגישה הלוקאלית שמנתחת את התקדמות הפרמטרים בכל צעד, לבין הגישה הגלובלית שמתארת את כל ההתפלגות, כזרימה מתמשכת של הסתברות במרחב המשקולות. בדיוק כמו בפיזיקה, המעבר מתיאור נקודתי לתיאור מבוזר מגלה תובנות שהיו נסתרות קודם. פתאום אפשר לשאול לא רק לאן המשקולות הולכות, אלא איפה הם מרוכזים, איך הם מתפזרים, ואיך המבנה של פונקציית הפסד משפיע על זה.'

# --- Paragraph 7 ---
$d.Paragraphs.Item(7).Range.Text = 'אחד החלקים המרשימים במאמר הוא הניתוח של רקורסיה בזמן. הכותבים לא מסתפקים בכך ש-SGD מתכנס, אלא בוחנים איך המבנה החוזר של תהליך הלמידה, המבוסס על חזרה עקבית דרך שיפועי הפונקציה, מתכתב עם הדינמיקה הרציפה של הפתרון למשוואות הפיזיקליות. דווקא ההשוואה הזו בין תהליך רקורסיבי עם צעד זמן דיסקרטי לבין תהליך דיפוזיה רציף מאפשרת לנסח לראשונה עקרונות כלליים על האפקטיביות של SGD: מתי הוא מצליח, מתי הוא עלול לסטות, וכיצד ניתן לשלוט בזה.'

# --- Paragraph 8 ---
$d.Paragraphs.Item(8).Range.Text = 'אבל מה שהכי תפס אותי הוא שהתמונה הזאת פותחת דלת לפיתוח עתידי. אם מקבלים את הפרדיגמה ש-SGD הוא לא רק תהליך חמדני שנע כלפי מטה, אלא מערכת פיזיקלית שמתפתחת לפי חוקים דיפרנציאליים אפשר להתחיל לתכנן אופטימיזציות חדשות מתוך אותו עולם מושגים. אולי לא צריך לשפר את SGD כמו שהוא, אלא לעבור ל-PDE-guided training, שבו מתארים ישירות את האבולוציה הרצויה של ההתפלגות, ופותרים אחורה כדי למצוא את הדינמיקה.'

# --- Paragraph 9 ---
$d.Paragraphs.Item(9).Range.Text = 'במובן הזה, המאמר הזה לא רק מסביר את העבר של SGD, אלא מציע עתיד חדש ללמידה עמוקה. עתיד שבו אנחנו פחות מגששים בתוך משטחים מרובי מימדים, ויותר בונים מודלים דינמיים עם מבנה פיזיקלי מובהק. זהו לא פחות שינוי תודעתי (ואולי גם פרקטי) שיכול לשנות את הדרך שבה ניגשים לאופטימיזציה כולה.'

# --- Remove the now-obsolete trailing paragraphs ---
# Results section through the concluding paragraph under the "מסקנה" heading
# (originally paragraphs 10-17, 1-indexed) is deleted entirely; repeatedly
# deleting at the same index removes each one in turn.
for ($i = 0; $i -lt 8; $i++) {
    $d.Paragraphs.Item(10).Range.Delete() | Out-Null
}

# --- Final paragraph: arXiv link ---
$d.Paragraphs.Item(10).Range.Text = 'https://arxiv.org/abs/2501.08425'
